$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-25 10:56:38"
$wsZhCn.Range("G3").Value = "2016-01-25 10:57:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-25 10:56:48"
$wsDeDe.Range("G3").Value = "2016-01-25 10:57:37"
